{"js": "const pairs = [\n  [\"60\u00d765=\", \"31\u00d759=\"],\n  [\"95\u00d781=\", \"53\u00d748=\"],\n  [\"50\u00d729=\", \"92\u00d752=\"],\n  [\"98\u00d778=\", \"65\u00d721=\"],\n  [\"28\u00d738=\", \"33\u00d769=\"],\n  [\"46\u00d733=\", \"70\u00d792=\"],\n  [\"32\u00d713=\", \"33\u00d736=\"],\n  [\"64\u00d760=\", \"79\u00d746=\"],\n  [\"90\u00d743=\", \"37\u00d776=\"],\n  [\"28\u00d758=\", \"50\u00d781=\"],\n  [\"32\u00d741=\", \"35\u00d781=\"],\n  [\"30\u00d715=\", \"28\u00d794=\"],\n  [\"66\u00d762=\", \"22\u00d764=\"],\n  [\"12\u00d792=\", \"31\u00d760=\"],\n  [\"76\u00d796=\", \"50\u00d777=\"],\n  [\"60\u00d759=\", \"95\u00d757=\"],\n  [\"32\u00d793=\", \"92\u00d733=\"],\n  [\"97\u00d738=\", \"12\u00d716=\"],\n  [\"83\u00d744=\", \"80\u00d712=\"],\n  [\"73\u00d723=\", \"12\u00d726=\"],\n  [\"49\u00d792=\", \"79\u00d720=\"],\n  [\"35\u00d781=\", \"48\u00d719=\"],\n  [\"99\u00d771=\", \"78\u00d739=\"],\n  [\"20\u00d772=\", \"85\u00d798=\"],\n  [\"22\u00d747=\", \"42\u00d765=\"],\n];\n\nconst body = context.document.body;\n\n// Collect the search-result range for each old value first (all old values\n// are unique within the document, so a single exact-match search per pair\n// is unambiguous even though some new values equal other pairs' old values).\nconst searches = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearches.forEach((s) => s.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [, newText] = pairs[i];\n  const results = searches[i];\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"AxB=\" multiplication prompts in the practice table with\n# their new values, preserving each run's existing formatting.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n  \"31\u00d759=\", \"53\u00d748=\", \"92\u00d752=\", \"65\u00d721=\", \"33\u00d769=\",\n  \"70\u00d792=\", \"33\u00d736=\", \"79\u00d746=\", \"37\u00d776=\", \"50\u00d781=\",\n  \"35\u00d781=\", \"28\u00d794=\", \"22\u00d764=\", \"31\u00d760=\", \"50\u00d777=\",\n  \"95\u00d757=\", \"92\u00d733=\", \"12\u00d716=\", \"80\u00d712=\", \"12\u00d726=\",\n  \"79\u00d720=\", \"48\u00d719=\", \"78\u00d739=\", \"85\u00d798=\", \"42\u00d765=\"\n)\n\n$dataRows = @(1, 5, 10, 15, 20)\n\n$i = 0\nforeach ($row in $dataRows) {\n  for ($col = 1; $col -le 5; $col++) {\n    $cell = $tbl.Cell($row, $col)\n    $range = $cell.Range\n    # Trim the trailing cell-end mark so we only overwrite the visible text,\n    # which keeps the existing run formatting (font/size) intact.\n    $range.MoveEnd(1, -1) | Out-Null\n    $range.Text = $newValues[$i]\n    $i++\n  }\n}\n\nWrite-Output \"updated $i cells\"\n"}
